$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "BURDEOS CIUDAD LA SALLE"
$ws.Range("B2").Value = 46511
$ws.Range("C2").Value = 46532

$ws.Range("A3").Value = "CHAMONIX CIUDAD LA SALLE"
$ws.Range("B3").Value = 46219
$ws.Range("C3").Value = 46420

$ws.Range("A4").Value = "LA SCALA"
$ws.Range("B4").Value = 46688
$ws.Range("C4").Value = 46792

$ws.Columns.Item(1).ColumnWidth = 23.77734375

$ws.Range("E5").Select()
